# Update stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B24").Value = 6416
$ws.Range("C24").Value = 1006
$ws.Range("D24").Value = 5982729
$ws.Range("E24").Value = 932.4702306733167
$ws.Range("F24").Value = 9.376065461984307
$ws.Range("G24").Value = 4.248704663212433
$ws.Range("H24").Value = 26.73737419956734
